$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    11983600.40858057,
    24109499.24552681,
    36385240.74239202,
    48811083.22630218,
    61389405.94803565,
    74123088.03592029,
    87014375.78095023,
    100066676.7130346,
    113281744.0853664,
    126662959.383555,
    140212550.7232948,
    153933252.6521844,
    167827875.1414022,
    181898939.5941035,
    196148841.3668919,
    210580676.7268349,
    225196740.6073617,
    240000036.1786282,
    254992338.7524262,
    270177256.3732307,
    285556833.5977859,
    301133497.876685,
    316910073.7869202,
    332888704.1142461,
    349072373.0651312,
    365462746.8894116,
    382062875.9694748,
    398874825.7032524,
    415900522.2335806,
    433142630.9041706,
    450603188.8640594,
    468284073.4616466,
    486187584.4480807,
    504315531.6374006,
    522669695.9944144,
    541252122.7816453,
    560064301.142359,
    579107792.2060256,
    598384542.8456945,
    617895193.6290576,
    637641639.5571072,
    657624387.4790373,
    677845303.7044461,
    698304754.9861217,
    719003158.3312334,
    739942050.9946457,
    761121242.5934942,
    782541330.4552808,
    804202139.6493908,
    826104397.1788691,
    848247096.6515927,
    870630689.0828905,
    893254134.1860774,
    916116999.6715529,
    939218091.7780775,
    962428420.5297517,
    985873093.9075214,
    1009551360.54368,
    1033461190.850898,
    1057601172.63974,
    1081968964.4237,
    1106562352.087063,
    1131378578.620336,
    1156415296.086895,
    1181669283.883677,
    1207137246.440516,
    1232816368.899849,
    1258702594.797298,
    1284792030.589759,
    1311080961.415303,
    1333102775.749794,
    1355282352.104061,
    1377614045.435207,
    1400121480.823915,
    1422772561.567305,
    1445561702.839376,
    1468484040.404572,
    1491534612.178619,
    1514707898.52471,
    1537998045.590745,
    1554811347.614815,
    1571692128.284739,
    1588636031.96566,
    1605636991.306359,
    1622690444.289208,
    1639790873.234804,
    1656879441.007428,
    1673377293.397635,
    1687799094.619394,
    1698394489.691973,
    1723086059.783255,
    1742082021.866847,
    1755660937.550134,
    1764630981.510189,
    1770049701.388512,
    1773014987.429447,
    1774531400.710859,
    1775445855.607762,
    1776428028.995756,
    1777976707.031778,
    1833971289.572975,
    1888680040.152385,
    1942408564.249466,
    1995411704.122757,
    2047896215.955525,
    2100027379.125194,
    2151934727.020473,
    2203717868.132726,
    2255451357.50292,
    2307189739.50681,
    2334038951.169407,
    2361328881.139693,
    2389131467.005021,
    2417407254.395346,
    2446116160.784106,
    2475217986.766279,
    2504380766.443476,
    2530311946.729125,
    2544603888.71852,
    2537365108.506991,
    2502012062.029463,
    2437073211.189877,
    2345418915.0537,
    2232519259.345423,
    2104798971.362617,
    1968495490.975441,
    1829017074.474315,
    1690685801.24449,
    1556724665.371577,
    1429367296.066275,
    1310025229.908217,
    1199464382.860801,
    1097961744.556799,
    1005442795.341471,
    921595009.3949554,
    845948868.4646014,
    777944235.1159639,
    716978813.6702818,
    662439328.2265866,
    613725967.1601751,
    570265354.2269646,
    531519590.8534323,
    496994980.4737478,
    466236723.9006115,
    438833293.3217613,
    414414515.7395141,
    392649513.497398,
    373241819.6849245,
    355929678.1018269,
    75395366.42377806
)

$startRow = 303
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $values[$i]
}
